$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "61.973.00"
Set-TextValue $ws.Range("E2") "  +0.78%  "
Set-TextValue $ws.Range("D3") "3.414.84"
Set-TextValue $ws.Range("E3") "  +1.00%  "
Set-TextValue $ws.Range("E4") "  -0.12%  "
Set-TextValue $ws.Range("D5") "410.30"
Set-TextValue $ws.Range("E5") "  +0.73%  "
Set-TextValue $ws.Range("D6") "128.61"
Set-TextValue $ws.Range("E6") "  -4.79%  "
Set-TextValue $ws.Range("D7") "0.624"
Set-TextValue $ws.Range("E7") "  +4.83%  "
Set-TextValue $ws.Range("D8") "0.999"
Set-TextValue $ws.Range("E8") "  -0.15%  "
Set-TextValue $ws.Range("D9") "0.749"
Set-TextValue $ws.Range("E9") "  +11.23%  "
Set-TextValue $ws.Range("E10") "  +15.58%  "
Set-TextValue $ws.Range("D11") "42.85"
Set-TextValue $ws.Range("E11") "  +0.41%  "
Set-TextValue $ws.Range("E12") "  -0.50%  "
Set-TextValue $ws.Range("D13") "21.26"
Set-TextValue $ws.Range("E13") "  +7.79%  "
Set-TextValue $ws.Range("D14") "8.84"
Set-TextValue $ws.Range("E14") "  +5.04%  "
Set-TextValue $ws.Range("D15") "0.0000201"
Set-TextValue $ws.Range("E15") "  +57.57%  "
Set-TextValue $ws.Range("D16") "3.406.46"
Set-TextValue $ws.Range("E16") "  +1.31%  "
Set-TextValue $ws.Range("E17") "  +14.81%  "
Set-TextValue $ws.Range("E18") "  +3.89%  "
Set-TextValue $ws.Range("D19") "61.909.16"
Set-TextValue $ws.Range("E19") "  +0.68%  "
Set-TextValue $ws.Range("D20") "406.41"
Set-TextValue $ws.Range("E20") "  +29.21%  "
Set-TextValue $ws.Range("D21") "90.65"
Set-TextValue $ws.Range("E21") "  +6.78%  "
Set-TextValue $ws.Range("D22") "3.19"
Set-TextValue $ws.Range("E22") "  -0.61%  "
Set-TextValue $ws.Range("D23") "13.43"
Set-TextValue $ws.Range("E23") "  +4.83%  "
Set-TextValue $ws.Range("E24") "  +3.09%  "
Set-TextValue $ws.Range("D25") "33.14"
Set-TextValue $ws.Range("E25") "  +12.12%  "
Set-TextValue $ws.Range("E26") "  -0.04%  "
Set-TextValue $ws.Range("D27") "8.54"
Set-TextValue $ws.Range("E27") "  +2.33%  "
Set-TextValue $ws.Range("E28") "  -0.10%  "
Set-TextValue $ws.Range("D29") "2.75"
Set-TextValue $ws.Range("E29") "  +7.18%  "
Set-TextValue $ws.Range("E30") "  +0.24%  "
Set-TextValue $ws.Range("D31") "0.172"
Set-TextValue $ws.Range("E31") "  -0.40%  "
Set-TextValue $ws.Range("D32") "43.87"
Set-TextValue $ws.Range("E32") "  +7.69%  "
Set-TextValue $ws.Range("D33") "11.79"
Set-TextValue $ws.Range("E33") "  +3.71%  "
Set-TextValue $ws.Range("E34") "  -0.05%  "
Set-TextValue $ws.Range("D35") "0.0499"
Set-TextValue $ws.Range("E35") "  +3.63%  "
Set-TextValue $ws.Range("D36") "52.64"
Set-TextValue $ws.Range("E36") "  +1.45%  "
Set-TextValue $ws.Range("D37") "0.998"
Set-TextValue $ws.Range("E37") "  -0.05%  "
Set-TextValue $ws.Range("D38") "3.38"
Set-TextValue $ws.Range("E38") "  -1.40%  "
Set-TextValue $ws.Range("D39") "2.90"
Set-TextValue $ws.Range("E39") "  -1.01%  "
Set-TextValue $ws.Range("E40") "  +5.78%  "
Set-TextValue $ws.Range("D41") "0.314"
Set-TextValue $ws.Range("E41") "  +6.93%  "
Set-TextValue $ws.Range("D42") "140.50"
Set-TextValue $ws.Range("E42") "  +1.41%  "
Set-TextValue $ws.Range("E43") "  -0.24%  "
Set-TextValue $ws.Range("D44") "4.02"
Set-TextValue $ws.Range("E44") "  -0.18%  "
Set-TextValue $ws.Range("D45") "2.37"
Set-TextValue $ws.Range("E45") "  +6.62%  "
Set-TextValue $ws.Range("D46") "16.76"
Set-TextValue $ws.Range("E46") "  -0.04%  "
Set-TextValue $ws.Range("D47") "21.74"
Set-TextValue $ws.Range("E47") "  +1.92%  "
Set-TextValue $ws.Range("D48") "2.106.98"
Set-TextValue $ws.Range("E48") "  -0.99%  "
Set-TextValue $ws.Range("B49") "BEAM"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextValue $ws.Range("D49") "0.0373"
Set-TextValue $ws.Range("E49") "  +8.79%  "
Set-TextValue $ws.Range("B50") "ThetaToken"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D50") "1.92"
Set-TextValue $ws.Range("E50") "  -3.01%  "
Set-TextValue $ws.Range("D51") "0.125"
Set-TextValue $ws.Range("E51") "  +12.42%  "
